$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 895.0833
$ws.Range("J19").Value = 546
$ws.Range("L19").Value = 546
$ws.Range("N19").Value = -896
$ws.Range("H33").Value = 126.4375
$ws.Range("I33").Value = 132.86667
$ws.Range("K33").Value = 132.86667
$ws.Range("M33").Value = 96.13333
$ws.Range("H74").Value = 5000
$ws.Range("I74").Value = 5000
$ws.Range("K74").Value = 5000
$ws.Range("M74").Value = -4064
$ws.Range("H77").Value = 5000
$ws.Range("I77").Value = 5000
$ws.Range("K77").Value = 25000
$ws.Range("M77").Value = -20320
$ws.Range("H137").Value = 2540.36
$ws.Range("I137").Value = 953.26666
$ws.Range("J137").Value = 4921
$ws.Range("K137").Value = 2859.79998
$ws.Range("L137").Value = 14763
$ws.Range("M137").Value = -309.7999799999998
$ws.Range("N137").Value = -19863
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3001.05
$ws.Range("I61").Value = 2889.4443
$ws.Range("K61").Value = 2889.4443
$ws.Range("M61").Value = -2677.4443
$ws.Range("H63").Value = 4143.909
$ws.Range("I63").Value = 3697.125
$ws.Range("J63").Value = 5335.3335
$ws.Range("K63").Value = 3697.125
$ws.Range("L63").Value = 5335.3335
$ws.Range("M63").Value = -3011.125
$ws.Range("N63").Value = -6707.3335
$ws.Range("H66").Value = 4143.909
$ws.Range("I66").Value = 3697.125
$ws.Range("J66").Value = 5335.3335
$ws.Range("K66").Value = 18485.625
$ws.Range("L66").Value = 26676.6675
$ws.Range("M66").Value = -15053.625
$ws.Range("N66").Value = -33540.6675
$ws.Range("H74").Value = 1514.875
$ws.Range("I74").Value = 1158.96
$ws.Range("K74").Value = 1158.96
$ws.Range("M74").Value = -284.96
$ws.Range("H77").Value = 1514.875
$ws.Range("I77").Value = 1158.96
$ws.Range("K77").Value = 5794.8
$ws.Range("M77").Value = -1426.8
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H122").Value = 668855
$ws.Range("I122").Value = 1112591.6
$ws.Range("K122").Value = 3337774.8
$ws.Range("M122").Value = -3335324.8
$ws.Range("H136").Value = 3001.05
$ws.Range("I136").Value = 2889.4443
$ws.Range("K136").Value = 8668.332900000001
$ws.Range("M136").Value = -6118.332900000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4915.4443
$ws.Range("J86").Value = 5360.8
$ws.Range("L86").Value = 5360.8
$ws.Range("N86").Value = -7606.8
$ws.Range("H89").Value = 4915.4443
$ws.Range("J89").Value = 5360.8
$ws.Range("L89").Value = 26804
$ws.Range("N89").Value = -38036
$ws.Range("H99").Value = 2491.5667
$ws.Range("I99").Value = 2978.5
$ws.Range("K99").Value = 2978.5
$ws.Range("M99").Value = -1480.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4615.6875
$ws.Range("I31").Value = 1464
$ws.Range("J31").Value = 5666.25
$ws.Range("K31").Value = 1464
$ws.Range("L31").Value = 5666.25
$ws.Range("M31").Value = -1169
$ws.Range("N31").Value = -6256.25
$ws.Range("H34").Value = 4615.6875
$ws.Range("I34").Value = 1464
$ws.Range("J34").Value = 5666.25
$ws.Range("K34").Value = 1464
$ws.Range("L34").Value = 5666.25
$ws.Range("M34").Value = -1262
$ws.Range("N34").Value = -6070.25
$ws.Range("H58").Value = 3201.4375
$ws.Range("I58").Value = 1998.75
$ws.Range("J58").Value = 3602.3333
$ws.Range("K58").Value = 1998.75
$ws.Range("L58").Value = 3602.3333
$ws.Range("M58").Value = -1795.75
$ws.Range("N58").Value = -4008.3333
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
$ws.Range("H136").Value = 3201.4375
$ws.Range("I136").Value = 1998.75
$ws.Range("J136").Value = 3602.3333
$ws.Range("K136").Value = 5996.25
$ws.Range("L136").Value = 10806.9999
$ws.Range("M136").Value = -3446.25
$ws.Range("N136").Value = -15906.9999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1248.1875
$ws.Range("J33").Value = 772.3
$ws.Range("L33").Value = 4633.799999999999
$ws.Range("N33").Value = -5199.799999999999
$ws.Range("H38").Value = 110.71429
$ws.Range("I38").Value = 71.75
$ws.Range("J38").Value = 162.66667
$ws.Range("K38").Value = 215.25
$ws.Range("L38").Value = 488.00001
$ws.Range("M38").Value = 131.75
$ws.Range("N38").Value = -1182.00001
$ws.Range("H70").Value = 2349
$ws.Range("I70").Value = 2018.8
$ws.Range("K70").Value = 6056.4
$ws.Range("M70").Value = -5741.4
$ws.Range("H73").Value = 2349
$ws.Range("I73").Value = 2018.8
$ws.Range("K73").Value = 6056.4
$ws.Range("M73").Value = -4964.4
$ws.Range("H76").Value = 4500
$ws.Range("J76").Value = 4500
$ws.Range("L76").Value = 13500
$ws.Range("N76").Value = -14266
$ws.Range("H79").Value = 4500
$ws.Range("J79").Value = 4500
$ws.Range("L79").Value = 13500
$ws.Range("N79").Value = -16152
$ws.Range("H80").Value = 2200
$ws.Range("I80").Value = 2200
$ws.Range("K80").Value = 6600
$ws.Range("M80").Value = -5664
$ws.Range("H83").Value = 2200
$ws.Range("I83").Value = 2200
$ws.Range("K83").Value = 19800
$ws.Range("M83").Value = -15120
$ws.Range("H121").Value = 789.75
$ws.Range("I121").Value = 500
$ws.Range("K121").Value = 1500
$ws.Range("M121").Value = -190
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2345.9429
$ws.Range("I132").Value = 1544.8422
$ws.Range("J132").Value = 3297.25
$ws.Range("K132").Value = 4634.5266
$ws.Range("L132").Value = 9891.75
$ws.Range("M132").Value = -2104.5266
$ws.Range("N132").Value = -14951.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1050.8572
$ws.Range("I22").Value = 713.6667
$ws.Range("J22").Value = 1657.8
$ws.Range("K22").Value = 713.6667
$ws.Range("L22").Value = 1657.8
$ws.Range("M22").Value = -418.6667
$ws.Range("N22").Value = -2247.8
$ws.Range("H27").Value = 1050.8572
$ws.Range("I27").Value = 713.6667
$ws.Range("J27").Value = 1657.8
$ws.Range("K27").Value = 713.6667
$ws.Range("L27").Value = 1657.8
$ws.Range("M27").Value = -606.6667
$ws.Range("N27").Value = -1871.8
$ws.Range("H55").Value = 412.73685
$ws.Range("J55").Value = 900
$ws.Range("L55").Value = 900
$ws.Range("N55").Value = -1246
$ws.Range("H122").Value = 5191.2666
$ws.Range("I122").Value = 4749.4
$ws.Range("K122").Value = 14248.2
$ws.Range("M122").Value = -11798.2
$ws.Range("H132").Value = 6066.6665
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 514500
$ws.Range("I3").Value = 1000000
$ws.Range("K3").Value = 1000000
$ws.Range("M3").Value = -999886
$ws.Range("H122").Value = 2989.1177
$ws.Range("I122").Value = 3117.7334
$ws.Range("K122").Value = 9353.200199999999
$ws.Range("M122").Value = -6903.200199999999
$ws.Range("H126").Value = 4093.3333
$ws.Range("I126").Value = 4081.3635
$ws.Range("J126").Value = 4126.25
$ws.Range("K126").Value = 12244.0905
$ws.Range("L126").Value = 12378.75
$ws.Range("M126").Value = -9774.0905
$ws.Range("N126").Value = -17318.75
$ws.Range("H132").Value = 1560.6
$ws.Range("I132").Value = 1449.5
$ws.Range("K132").Value = 4348.5
$ws.Range("M132").Value = -1818.5
